$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date labels in column A (rows 2-7): force text so Excel
# does not auto-convert these ISO date strings into date serials.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-09-30"
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-10-15"
$ws.Range("A3").Style = "Normal"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-10-30"
$ws.Range("A4").Style = "Normal"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-11-15"
$ws.Range("A5").Style = "Normal"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-11-30"
$ws.Range("A6").Style = "Normal"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-12-15"
$ws.Range("A7").Style = "Normal"

# Update the numeric grid (columns B:K, rows 2-7) with the new nowcast values.
$ws.Range("B2").Value = [double]"0.27542368938497669"
$ws.Range("C2").Value = [double]"0"
$ws.Range("D2").Value = [double]"0"
$ws.Range("E2").Value = [double]"0"
$ws.Range("F2").Value = [double]"0"
$ws.Range("G2").Value = [double]"0"
$ws.Range("H2").Value = [double]"0"
$ws.Range("I2").Value = [double]"0"
$ws.Range("J2").Value = [double]"0"
$ws.Range("K2").Value = [double]"0"

$ws.Range("B3").Value = [double]"0.27600624442285482"
$ws.Range("C3").Value = [double]"0"
$ws.Range("D3").Value = [double]"-0.00054505057472851006"
$ws.Range("E3").Value = [double]"0.00095968984101283945"
$ws.Range("F3").Value = [double]"0.0016516213937664231"
$ws.Range("G3").Value = [double]"0.001387991156290887"
$ws.Range("H3").Value = [double]"-4.0940504011610656e-05"
$ws.Range("I3").Value = [double]"-0.00017837202157341476"
$ws.Range("J3").Value = [double]"0"
$ws.Range("K3").Value = [double]"0.00012049809947778378"

$ws.Range("B4").Value = [double]"0.35721374530323058"
$ws.Range("C4").Value = [double]"0.056517975558058542"
$ws.Range("D4").Value = [double]"0"
$ws.Range("E4").Value = [double]"-0.00019298174139811754"
$ws.Range("F4").Value = [double]"2.1033098067586637e-06"
$ws.Range("G4").Value = [double]"0"
$ws.Range("H4").Value = [double]"0.00064808978942621919"
$ws.Range("I4").Value = [double]"0.0017505463688082328"
$ws.Range("J4").Value = [double]"0.0033534147901798732"
$ws.Range("K4").Value = [double]"-6.5464596027742239e-05"

$ws.Range("B5").Value = [double]"0.30733777652680849"
$ws.Range("C5").Value = [double]"0"
$ws.Range("D5").Value = [double]"0.00034843227590949406"
$ws.Range("E5").Value = [double]"-0.0016054093380629387"
$ws.Range("F5").Value = [double]"0.011643809053714666"
$ws.Range("G5").Value = [double]"-0.0030629124637942345"
$ws.Range("H5").Value = [double]"0.00021054422943723685"
$ws.Range("I5").Value = [double]"-0.00027013226793090372"
$ws.Range("J5").Value = [double]"0"
$ws.Range("K5").Value = [double]"-0.011989028776616228"

$ws.Range("B6").Value = [double]"0.22535638213167042"
$ws.Range("C6").Value = [double]"-0.029583831479153687"
$ws.Range("D6").Value = [double]"0"
$ws.Range("E6").Value = [double]"-0.0046189933715279846"
$ws.Range("F6").Value = [double]"0.00066098148898444619"
$ws.Range("G6").Value = [double]"0"
$ws.Range("H6").Value = [double]"0.0016721532257435257"
$ws.Range("I6").Value = [double]"-0.013921219016512229"
$ws.Range("J6").Value = [double]"0"
$ws.Range("K6").Value = [double]"-1.1684333257322432e-05"

$ws.Range("B7").Value = [double]"0.25866158402573242"
$ws.Range("C7").Value = [double]"0"
$ws.Range("D7").Value = [double]"-0.033165550624100273"
$ws.Range("E7").Value = [double]"-0.00046343405080691372"
$ws.Range("F7").Value = [double]"0.0096096374886073047"
$ws.Range("G7").Value = [double]"0.0021870738249456386"
$ws.Range("H7").Value = [double]"0"
$ws.Range("I7").Value = [double]"0"
$ws.Range("J7").Value = [double]"0"
$ws.Range("K7").Value = [double]"-0.0010245163857025008"

# Column widths were re-auto-fitted for the refreshed content.
$ws.Columns.Item(6).ColumnWidth = 14.333333333333332
$ws.Columns.Item(7).ColumnWidth = 14.333333333333332
$ws.Columns.Item(8).ColumnWidth = 15.0
$ws.Columns.Item(10).ColumnWidth = 14.166666666666668
$ws.Columns.Item(11).ColumnWidth = 15.0
